$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-04-06 Saturday" "2024-04-07 Sunday"

Replace-Text "336×3=" "651×7="
Replace-Text "570×6=" "624×7="
Replace-Text "674×3=" "144×3="
Replace-Text "854×5=" "291×8="
Replace-Text "278×9=" "675×6="
Replace-Text "954×7=" "321×6="
Replace-Text "327×2=" "200×4="
Replace-Text "564×6=" "389×4="
Replace-Text "123×3=" "302×8="
Replace-Text "442×6=" "873×3="
Replace-Text "165×8=" "584×5="
Replace-Text "785×5=" "838×5="
Replace-Text "989×3=" "106×4="
Replace-Text "361×8=" "634×5="
Replace-Text "635×9=" "464×5="
Replace-Text "831×8=" "648×6="
Replace-Text "264×7=" "515×3="
Replace-Text "678×9=" "465×6="
Replace-Text "301×6=" "144×8="
Replace-Text "963×2=" "904×4="
Replace-Text "246×4=" "178×7="
Replace-Text "419×7=" "427×4="
Replace-Text "553×8=" "742×7="
Replace-Text "408×3=" "488×4="
Replace-Text "680×6=" "401×8="

Write-Output "Done replacing all values"
